$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 138 (IR_MLID) to hold the new
# "MonitoringLocationIdentifier" translation entry, shifting the rest of
# the table (rows 138-164) down by one.
$ws.Rows.Item(138).Insert()

# Populate the new row: column A gets the new source-column name, column B
# keeps the same "DS" sheet tag used by the surrounding rows.
$ws.Range("A138").Value = "MonitoringLocationIdentifier"
$ws.Range("B138").Value = "DS"

# Reflect the scrolled/selected state seen after the edit.
[void]$ws.Activate()
[void]$ws.Range("B139").Select()
$excel.ActiveWindow.ScrollRow = 125
$excel.ActiveWindow.ScrollColumn = 1
